$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new solar system card rows, appended after the existing 47 data rows
# (rows 2-47), mirroring the layout/style of the other "solar" type rows.
$newRows = @(
    @{ idx = 47; title = "Appert System";   desc = "Do not mess with us";                    image = "Appert System.png" },
    @{ idx = 48; title = "Cartof System";    desc = "We only like peaceful people here";      image = "Cartof System.png" },
    @{ idx = 49; title = "VanWeiss System";  desc = "Almost time travels through space";      image = "VanWeiss System.png" }
)

$row = 48
foreach ($card in $newRows) {
    $ws.Cells.Item($row, 1).Value = $card.idx
    $ws.Cells.Item($row, 2).Value = "solar"
    $ws.Cells.Item($row, 3).Value = "null"
    $ws.Cells.Item($row, 4).Value = $card.title
    $ws.Cells.Item($row, 6).Value = 5
    $ws.Cells.Item($row, 7).Value = "null"
    $ws.Cells.Item($row, 8).Value = 0.03
    $ws.Cells.Item($row, 9).Value = $card.image
    $ws.Cells.Item($row, 10).Value = 100

    $row = $row + 1
}

$row = 48
foreach ($card in $newRows) {
    $ws.Cells.Item($row, 5).Value = $card.desc
    $row = $row + 1
}

# Match the author's final view state: scrolled down with J51 selected.
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J51").Select()
